$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'326.09"
$ws.Range("E2").Value = "'-1.18%"
$ws.Range("G2").Value = "'11"

# Row 3
$ws.Range("D3").Value = "'39.62"
$ws.Range("E3").Value = "'-1.21%"
$ws.Range("G3").Value = "'11"

# Row 4
$ws.Range("D4").Value = "'5.721"
$ws.Range("E4").Value = "'6.63%"
$ws.Range("G4").Value = "'11"

# Row 5
$ws.Range("D5").Value = "'0.08037"
$ws.Range("E5").Value = "'-0.87%"
$ws.Range("G5").Value = "'11"

# Row 6
$ws.Range("D6").Value = "'2.054"
$ws.Range("E6").Value = "'6.85%"
$ws.Range("G6").Value = "'11"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.492"
$ws.Range("E7").Value = "'-0.79%"
$ws.Range("G7").Value = "'11"

# Row 8
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.631"
$ws.Range("E8").Value = "'-0.22%"
$ws.Range("G8").Value = "'11"

# Row 9
$ws.Range("D9").Value = "'2.971"
$ws.Range("E9").Value = "'-0.22%"
$ws.Range("G9").Value = "'11"

# Row 10
$ws.Range("D10").Value = "'0.9221"
$ws.Range("E10").Value = "'-2.22%"
$ws.Range("G10").Value = "'11"

# Row 11
$ws.Range("D11").Value = "'0.1247"
$ws.Range("E11").Value = "'-8.39%"
$ws.Range("G11").Value = "'11"

# Row 12
$ws.Range("D12").Value = "'0.1956"
$ws.Range("E12").Value = "'-1.02%"
$ws.Range("G12").Value = "'11"

# Row 13
$ws.Range("D13").Value = "'8.739"
$ws.Range("E13").Value = "'21.06%"
$ws.Range("G13").Value = "'11"

# Row 14
$ws.Range("D14").Value = "'0.09187"
$ws.Range("E14").Value = "'-1.61%"
$ws.Range("G14").Value = "'11"

# Row 15
$ws.Range("D15").Value = "'0.03568"
$ws.Range("E15").Value = "'0.19%"
$ws.Range("G15").Value = "'11"

# Row 16
$ws.Range("D16").Value = "'0.1049"
$ws.Range("E16").Value = "'9.36%"
$ws.Range("G16").Value = "'11"

# Row 17
$ws.Range("D17").Value = "'0.001296"
$ws.Range("E17").Value = "'-2.08%"
$ws.Range("G17").Value = "'11"

# Row 18
$ws.Range("D18").Value = "'0.006122"
$ws.Range("E18").Value = "'-3.62%"
$ws.Range("G18").Value = "'11"

# Row 19
$ws.Range("D19").Value = "'3.350"
$ws.Range("E19").Value = "'-0.50%"
$ws.Range("G19").Value = "'11"

# Row 20
$ws.Range("E20").Value = "'-1.11%"
$ws.Range("G20").Value = "'11"

# Row 21
$ws.Range("D21").Value = "'0.1351"
$ws.Range("E21").Value = "'1.43%"
$ws.Range("G21").Value = "'11"

# Row 22
$ws.Range("D22").Value = "'0.2389"
$ws.Range("E22").Value = "'-6.71%"
$ws.Range("G22").Value = "'11"

# Row 23
$ws.Range("D23").Value = "'0.04390"
$ws.Range("E23").Value = "'-0.82%"
$ws.Range("G23").Value = "'11"

# Row 24
$ws.Range("D24").Value = "'0.001259"
$ws.Range("E24").Value = "'3.15%"
$ws.Range("G24").Value = "'11"

# Row 25
$ws.Range("D25").Value = "'0.004607"
$ws.Range("E25").Value = "'7.67%"
$ws.Range("G25").Value = "'11"

# Row 26
$ws.Range("E26").Value = "'2.48%"
$ws.Range("G26").Value = "'11"

# Row 27
$ws.Range("G27").Value = "'11"

# Row 28
$ws.Range("G28").Value = "'11"

# Row 29
$ws.Range("G29").Value = "'11"

# Row 30
$ws.Range("G30").Value = "'11"

# Row 31
$ws.Range("G31").Value = "'11"

# Row 32
$ws.Range("G32").Value = "'11"

# Row 33
$ws.Range("G33").Value = "'11"

# Row 34
$ws.Range("G34").Value = "'11"

# Row 35
$ws.Range("G35").Value = "'11"

# Row 36
$ws.Range("G36").Value = "'11"

# Row 37
$ws.Range("G37").Value = "'11"

# Row 38
$ws.Range("G38").Value = "'11"

# Row 39
$ws.Range("D39").Value = "'0.02501"
$ws.Range("E39").Value = "'0.71%"
$ws.Range("G39").Value = "'11"

# Row 40
$ws.Range("D40").Value = "'0.05318"
$ws.Range("E40").Value = "'1.90%"
$ws.Range("G40").Value = "'11"

# Row 41
$ws.Range("D41").Value = "'0.007480"
$ws.Range("E41").Value = "'-0.81%"
$ws.Range("G41").Value = "'11"

# Row 42
$ws.Range("D42").Value = "'0.009910"
$ws.Range("E42").Value = "'8.92%"
$ws.Range("G42").Value = "'11"

# Row 43
$ws.Range("D43").Value = "'0.1406"
$ws.Range("E43").Value = "'-1.67%"
$ws.Range("G43").Value = "'11"

# Row 44
$ws.Range("D44").Value = "'0.002117"
$ws.Range("E44").Value = "'-2.45%"
$ws.Range("G44").Value = "'11"

# Row 45
$ws.Range("D45").Value = "'0.01110"
$ws.Range("E45").Value = "'2.19%"
$ws.Range("G45").Value = "'11"

# Row 46
$ws.Range("D46").Value = "'0.00006694"
$ws.Range("E46").Value = "'1.12%"
$ws.Range("G46").Value = "'11"

# Row 47
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("G47").Value = "'11"

# Row 48
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.003039"
$ws.Range("E48").Value = "'-9.14%"
$ws.Range("G48").Value = "'11"

# Row 49
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.002279"
$ws.Range("E49").Value = "'-5.05%"
$ws.Range("G49").Value = "'11"

# Row 50
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("G50").Value = "'11"

# Row 51
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.07%"
$ws.Range("G51").Value = "'11"
